$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 861.8461
$ws.Range("I55").Value = 275.5
$ws.Range("K55").Value = 275.5
$ws.Range("M55").Value = -61.5
$ws.Range("H80").Value = 1475.5883
$ws.Range("I80").Value = 1052.7142
$ws.Range("K80").Value = 3158.1426
$ws.Range("M80").Value = -2160.1426
$ws.Range("H83").Value = 1475.5883
$ws.Range("I83").Value = 1052.7142
$ws.Range("K83").Value = 9474.427799999999
$ws.Range("M83").Value = -4482.427799999999
$ws.Range("H86").Value = 1542.8462
$ws.Range("I86").Value = 1428.7778
$ws.Range("J86").Value = 1799.5
$ws.Range("K86").Value = 1428.7778
$ws.Range("L86").Value = 1799.5
$ws.Range("M86").Value = -305.7778000000001
$ws.Range("N86").Value = -4045.5
$ws.Range("H88").Value = 2499.8572
$ws.Range("J88").Value = 1999.6666
$ws.Range("L88").Value = 1999.6666
$ws.Range("N88").Value = -2811.6666
$ws.Range("H89").Value = 1542.8462
$ws.Range("I89").Value = 1428.7778
$ws.Range("J89").Value = 1799.5
$ws.Range("K89").Value = 7143.889
$ws.Range("L89").Value = 8997.5
$ws.Range("M89").Value = -1527.889
$ws.Range("N89").Value = -20229.5
$ws.Range("H91").Value = 2499.8572
$ws.Range("J91").Value = 1999.6666
$ws.Range("L91").Value = 1999.6666
$ws.Range("N91").Value = -4807.6666
$ws.Range("H132").Value = 3185.1924
$ws.Range("I132").Value = 2848.4211
$ws.Range("K132").Value = 8545.263300000001
$ws.Range("M132").Value = -6015.263300000001
$ws.Range("H141").Value = 6659.2
$ws.Range("I141").Value = 7900.5
$ws.Range("K141").Value = 23701.5
$ws.Range("M141").Value = -18521.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 411
$ws.Range("I2").Value = 370.25
$ws.Range("K2").Value = 370.25
$ws.Range("M2").Value = -257.25
$ws.Range("H32").Value = 9195
$ws.Range("I32").Value = 8594.723
$ws.Range("K32").Value = 8594.723
$ws.Range("M32").Value = -8307.723
$ws.Range("H74").Value = 3756.353
$ws.Range("I74").Value = 2950
$ws.Range("J74").Value = 4908.2856
$ws.Range("K74").Value = 2950
$ws.Range("L74").Value = 4908.2856
$ws.Range("M74").Value = -2076
$ws.Range("N74").Value = -6656.2856
$ws.Range("H77").Value = 3756.353
$ws.Range("I77").Value = 2950
$ws.Range("J77").Value = 4908.2856
$ws.Range("K77").Value = 14750
$ws.Range("L77").Value = 24541.428
$ws.Range("M77").Value = -10382
$ws.Range("N77").Value = -33277.428
$ws.Range("H110").Value = 5980.276
$ws.Range("I110").Value = 6143.5
$ws.Range("K110").Value = 6143.5
$ws.Range("M110").Value = -4098.5
$ws.Range("H116").Value = 411
$ws.Range("I116").Value = 370.25
$ws.Range("K116").Value = 370.25
$ws.Range("M116").Value = 1923.75
$ws.Range("H122").Value = 3345
$ws.Range("I122").Value = 5149.1665
$ws.Range("J122").Value = 1180
$ws.Range("K122").Value = 15447.4995
$ws.Range("L122").Value = 3540
$ws.Range("M122").Value = -12997.4995
$ws.Range("N122").Value = -8440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 411
$ws.Range("I3").Value = 370.25
$ws.Range("K3").Value = 370.25
$ws.Range("M3").Value = -256.25
$ws.Range("H22").Value = 381.25
$ws.Range("J22").Value = 275
$ws.Range("L22").Value = 275
$ws.Range("N22").Value = -621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1932.8334
$ws.Range("I16").Value = 1719.4
$ws.Range("K16").Value = 1719.4
$ws.Range("M16").Value = -1432.4
$ws.Range("H22").Value = 4001898.2
$ws.Range("I22").Value = 1776.4
$ws.Range("J22").Value = 8002020
$ws.Range("K22").Value = 1776.4
$ws.Range("L22").Value = 8002020
$ws.Range("M22").Value = -1426.4
$ws.Range("N22").Value = -8002720
$ws.Range("H31").Value = 3157.1667
$ws.Range("I31").Value = 2577.5557
$ws.Range("K31").Value = 2577.5557
$ws.Range("M31").Value = -2282.5557
$ws.Range("H34").Value = 3157.1667
$ws.Range("I34").Value = 2577.5557
$ws.Range("K34").Value = 2577.5557
$ws.Range("M34").Value = -2375.5557
$ws.Range("H58").Value = 2698.1428
$ws.Range("I58").Value = 1481.1666
$ws.Range("K58").Value = 1481.1666
$ws.Range("M58").Value = -1278.1666
$ws.Range("H63").Value = 100271
$ws.Range("J63").Value = 100271
$ws.Range("L63").Value = 100271
$ws.Range("N63").Value = -101643
$ws.Range("H66").Value = 100271
$ws.Range("J66").Value = 100271
$ws.Range("L66").Value = 300813
$ws.Range("N66").Value = -307677
$ws.Range("H86").Value = 8133
$ws.Range("I86").Value = 8133
$ws.Range("K86").Value = 8133
$ws.Range("M86").Value = -7010
$ws.Range("H89").Value = 8133
$ws.Range("I89").Value = 8133
$ws.Range("K89").Value = 40665
$ws.Range("M89").Value = -35049
$ws.Range("H107").Value = 673.2778
$ws.Range("I107").Value = 481.66666
$ws.Range("K107").Value = 481.66666
$ws.Range("M107").Value = 1438.33334
$ws.Range("H113").Value = 1932.8334
$ws.Range("I113").Value = 1719.4
$ws.Range("K113").Value = 1719.4
$ws.Range("M113").Value = 450.5999999999999
$ws.Range("H132").Value = 2531.2
$ws.Range("I132").Value = 2283.4211
$ws.Range("J132").Value = 3315.8333
$ws.Range("K132").Value = 6850.263300000001
$ws.Range("L132").Value = 9947.499899999999
$ws.Range("M132").Value = -4320.263300000001
$ws.Range("N132").Value = -15007.4999
$ws.Range("H134").Value = 2542.6
$ws.Range("I134").Value = 2188.1667
$ws.Range("J134").Value = 3074.25
$ws.Range("K134").Value = 6564.500100000001
$ws.Range("L134").Value = 9222.75
$ws.Range("M134").Value = -4029.500100000001
$ws.Range("N134").Value = -14292.75
$ws.Range("H136").Value = 2698.1428
$ws.Range("I136").Value = 1481.1666
$ws.Range("K136").Value = 4443.4998
$ws.Range("M136").Value = -1893.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 650.5
$ws.Range("I5").Value = 659
$ws.Range("K5").Value = 1977
$ws.Range("M5").Value = -1865
$ws.Range("H11").Value = 96.545456
$ws.Range("I11").Value = 135.66667
$ws.Range("K11").Value = 407.00001
$ws.Range("M11").Value = -267.00001
$ws.Range("H135").Value = 650.5
$ws.Range("I135").Value = 659
$ws.Range("K135").Value = 5931
$ws.Range("M135").Value = -3396

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1840
$ws.Range("J122").Value = 1298.5
$ws.Range("L122").Value = 3895.5
$ws.Range("N122").Value = -8795.5
$ws.Range("H132").Value = 2627.2964
$ws.Range("I132").Value = 2514.818
$ws.Range("J132").Value = 3122.2
$ws.Range("K132").Value = 7544.454000000001
$ws.Range("L132").Value = 9366.599999999999
$ws.Range("M132").Value = -5014.454000000001
$ws.Range("N132").Value = -14426.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1735.5454
$ws.Range("I7").Value = 1735.5454
$ws.Range("K7").Value = 1735.5454
$ws.Range("M7").Value = -1623.5454
$ws.Range("H16").Value = 1330.2
$ws.Range("I16").Value = 625.5
$ws.Range("K16").Value = 625.5
$ws.Range("M16").Value = -455.5
$ws.Range("H46").Value = 1452.3158
$ws.Range("I46").Value = 1090.5294
$ws.Range("K46").Value = 1090.5294
$ws.Range("M46").Value = -902.5293999999999
$ws.Range("H61").Value = 1709.5
$ws.Range("I61").Value = 1566.1111
$ws.Range("K61").Value = 1566.1111
$ws.Range("M61").Value = -1364.1111
$ws.Range("H82").Value = 1485.4
$ws.Range("I82").Value = 699.6667
$ws.Range("J82").Value = 2009.2222
$ws.Range("K82").Value = 699.6667
$ws.Range("L82").Value = 2009.2222
$ws.Range("M82").Value = -338.6667
$ws.Range("N82").Value = -2731.2222
$ws.Range("H85").Value = 1485.4
$ws.Range("I85").Value = 699.6667
$ws.Range("J85").Value = 2009.2222
$ws.Range("K85").Value = 699.6667
$ws.Range("L85").Value = 2009.2222
$ws.Range("M85").Value = 548.3333
$ws.Range("N85").Value = -4505.2222
$ws.Range("H113").Value = 1709.5
$ws.Range("I113").Value = 1566.1111
$ws.Range("K113").Value = 1566.1111
$ws.Range("M113").Value = 603.8888999999999
$ws.Range("H126").Value = 1735.5454
$ws.Range("I126").Value = 1735.5454
$ws.Range("K126").Value = 5206.6362
$ws.Range("M126").Value = -2736.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3208.3845
$ws.Range("I132").Value = 2570.9
$ws.Range("J132").Value = 5333.3335
$ws.Range("K132").Value = 7712.700000000001
$ws.Range("L132").Value = 16000.0005
$ws.Range("M132").Value = -5182.700000000001
$ws.Range("N132").Value = -21060.0005
